# Fruta / hortaliza, semanal
# Insert a new weekly record at row 44 ("Macroferia Regional de Talca" /
# Damasco / Dina / Primera, $/caja 18 kilos, Región de O'Higgins) and push
# the existing rows 44-52 down to 45-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 44, shifting rows 44:52 down
# to 45:53 (mirrors Excel's Rows("44:44").Insert Shift:=xlShiftDown).
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new weekly price record.
$ws.Range("A44").Value = 5
$ws.Range("B44").Value = "Macroferia Regional de Talca"
$ws.Range("C44").Value = "Maule"
$ws.Range("D44").Value = 44918
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100103
$ws.Range("H44").Value = "Frutos de hueso (carozo)"
$ws.Range("I44").Value = 100103003
$ws.Range("J44").Value = "Damasco"
$ws.Range("K44").Value = "Dina"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 18000
$ws.Range("O44").Value = 18000
$ws.Range("P44").Value = 18000
$ws.Range("Q44").Value = "$/caja 18 kilos"
$ws.Range("R44").Value = "Región de O'Higgins"
$ws.Range("S44").Value = 1000
$ws.Range("T44").Value = 18
